$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '27.008.42'
$ws.Range("E2").Value = '  -0.71%  '
Set-TextValue $ws.Range("D3") '1.622.37'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws.Range("D5") '213.50'
$ws.Range("E5").Value = '  -1.63%  '
Set-TextValue $ws.Range("D6") '0.509'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  -1.51%  '
Set-TextValue $ws.Range("D10") '19.95'
$ws.Range("E10").Value = '  -0.84%  '
Set-TextValue $ws.Range("D11") '0.0838'
$ws.Range("E11").Value = '  -1.40%  '
Set-TextValue $ws.Range("D12") '1.849.92'
$ws.Range("E12").Value = '  -0.88%  '
Set-TextValue $ws.Range("D13") '1.619.32'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("E15").Value = '  -1.33%  '
Set-TextValue $ws.Range("D16") '26.988.31'
$ws.Range("E16").Value = '  -0.70%  '
Set-TextValue $ws.Range("D17") '64.17'
$ws.Range("E17").Value = '  -3.60%  '
Set-TextValue $ws.Range("D18") '0.0₃0734'
$ws.Range("E18").Value = '  -0.57%  '
Set-TextValue $ws.Range("D19") '213.33'
$ws.Range("E19").Value = '  -1.86%  '
$ws.Range("E20").Value = '  +0.08%  '
Set-TextValue $ws.Range("D21") '6.81'
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("E22").Value = '  -2.30%  '
$ws.Range("E23").Value = '  -7.52%  '
$ws.Range("E24").Value = '  -2.11%  '
Set-TextValue $ws.Range("D25") '146.73'
$ws.Range("E25").Value = '  -0.59%  '
Set-TextValue $ws.Range("D26") '7.46'
$ws.Range("E26").Value = '  +0.95%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -4.18%  '
Set-TextValue $ws.Range("D29") '15.49'
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("E32").Value = '  -2.74%  '
Set-TextValue $ws.Range("D33") '0.699'
$ws.Range("E33").Value = '  +27.28%  '
$ws.Range("E34").Value = '  -1.11%  '
Set-TextValue $ws.Range("D35") '1.339.04'
$ws.Range("E35").Value = '  +2.68%  '
Set-TextValue $ws.Range("D36") '1.55'
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("E38").Value = '  -0.79%  '
Set-TextValue $ws.Range("D39") '0.841'
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("E40").Value = '  -0.01%  '
Set-TextValue $ws.Range("D43") '5.36'
$ws.Range("E43").Value = '  +0.45%  '
Set-TextValue $ws.Range("D44") '63.72'
$ws.Range("E44").Value = '  +2.21%  '
Set-TextValue $ws.Range("D45") '1.760.77'
$ws.Range("E45").Value = '  -0.89%  '
Set-TextValue $ws.Range("D46") '89.78'
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("E47").Value = '  +1.49%  '
Set-TextValue $ws.Range("D48") '0.804'
$ws.Range("E48").Value = '  +8.88%  '
$ws.Range("E49").Value = '  +0.17%  '
Set-TextValue $ws.Range("D50") '0.0991'
$ws.Range("E50").Value = '  +3.50%  '
Set-TextValue $ws.Range("D51") '7.58'
$ws.Range("E51").Value = '  -0.92%  '

# Row 41 <-> Row 42 swap (TrustWalletToken <-> MXToken) plus independent E updates
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D41") '2.22'
$ws.Range("E41").Value = '  -0.91%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D42") '0.797'
$ws.Range("E42").Value = '  -1.11%  '
